$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 102, pushing existing rows 102:112 down to 103:113
$ws.Rows("102:102").Insert()

# Populate the newly inserted row 102 with the new record's data
$ws.Cells.Item(102, 1).Value = 11
$ws.Cells.Item(102, 2).Value = 'Vega Monumental Concepción'
$ws.Cells.Item(102, 3).Value = 'Bíobío'
$ws.Cells.Item(102, 4).Value = 44505
$ws.Cells.Item(102, 5).Value = 8
$ws.Cells.Item(102, 6).Value = 'Fruta'
$ws.Cells.Item(102, 7).Value = 100108
$ws.Cells.Item(102, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(102, 9).Value = 100108005
$ws.Cells.Item(102, 10).Value = 'Piña'
$ws.Cells.Item(102, 11).Value = 'Caramelo'
$ws.Cells.Item(102, 12).Value = 'Segunda'
$ws.Cells.Item(102, 13).Value = 200
$ws.Cells.Item(102, 14).Value = 17000
$ws.Cells.Item(102, 15).Value = 18000
$ws.Cells.Item(102, 16).Value = 17500
$ws.Cells.Item(102, 17).Value = '$/caja 14 unidades'
$ws.Cells.Item(102, 18).Value = 'Ecuador'
$ws.Cells.Item(102, 19).Value = 1250
$ws.Cells.Item(102, 20).Value = 14
